# Weekly update: insert a new "Acelga" price record (week of 2021-11-16)
# for Vega Monumental Concepcion, pushing the existing rows 167-178 down
# to 168-179.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 167, shifting rows
# 167:178 down to 168:179 (and the used range / dimension grows to R179).
$ws.Rows.Item(167).Insert()

# Populate the newly inserted row 167 with the new weekly record.
$ws.Range("A167").Value = 11
$ws.Range("B167").Value = "Vega Monumental Concepción"
$ws.Range("C167").Value = "Bíobío"
$ws.Range("D167").Value = 44516
$ws.Range("E167").Value = 8
$ws.Range("F167").Value = 100112009
$ws.Range("G167").Value = "Acelga"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 450
$ws.Range("K167").Value = 600
$ws.Range("L167").Value = 650
$ws.Range("M167").Value = 622
$ws.Range("N167").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O167").Value = "Región de Ñuble"
$ws.Range("P167").Value = 622
$ws.Range("Q167").Value = 1
$ws.Range("R167").Value = "Hortaliza"
